$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.936.50"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "1.856.37"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("D4").Value = "'1.002"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'317.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("D6").Value = "'1.001"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.4347"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  -5.27%  "
$ws.Range("D8").Value = "'0.3685"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D9").Value = "'0.07477"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").Value = "'0.9378"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "'21.27"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "1.831.49"
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("D13").Value = "'6.695"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("D14").Value = "'5.425"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  -4.12%  "
$ws.Range("D15").Value = "'0.06856"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'81.40"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "'0.000009038"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("D19").Value = "'1.001"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'15.86"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").Value = "27.908.02"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").Value = "'5.100"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "2.143.63"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'2.007"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("D26").Value = "'153.75"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").Value = "'18.35"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("D28").Value = "'5.383"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("D29").Value = "'113.34"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").Value = "'1.727"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -7.42%  "
$ws.Range("D31").Value = "'0.08959"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "'0.8025"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("D33").Value = "'4.829"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  -4.79%  "
$ws.Range("D34").Value = "'2.995"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("D36").Value = "'1.000"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'1.116"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "'0.05444"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  -4.70%  "
$ws.Range("D39").Value = "'0.01967"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("D40").Value = "'2.928"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "'0.5230"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").Value = "'6.993"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -5.75%  "
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("D44").Value = "'8.765"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  -5.92%  "
$ws.Range("D45").Value = "'0.06718"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("D47").Value = "'10.64"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").Value = "'106.11"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D49").Value = "'1.933"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  -7.35%  "
$ws.Range("D50").Value = "'1.675"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "'0.9995"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
